$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("E2E")

# Insert the 4 new rows top-down; each Insert() operates on the CURRENT sheet
# row at the time of the call, so using the final target row numbers in
# ascending order naturally produces the correct final layout.
$ws.Rows("44").Insert()
$ws.Rows("46").Insert()
$ws.Rows("49").Insert()
$ws.Rows("52").Insert()

# Populate cell values in the same chronological order the original author
# used, so new shared strings are appended to sst in the matching order.
$ws.Range("A49").Value = "SKU-11261400 - 2QTY"
$ws.Range("M49").Value = "OXO Good Grips Citrus Zester"
$ws.Range("N49").Value = "'2"
$ws.Range("AJ49").Value = "'2"
$ws.Range("AJ49").ClearContents()

$ws.Range("A52").Value = "SKU-32480 -2QTY"
$ws.Range("J52").Value = "'2"
$ws.Range("J52").ClearContents()
$ws.Range("M52").Value = "Salad Spinner"
$ws.Range("N52").Value = "'2"
$ws.Range("AJ52").Value = "'2"
$ws.Range("AJ52").ClearContents()

$ws.Range("A46").Value = "Bundle Product"
$ws.Range("M46").Value = "OXO 5-Piece Barware Set - Exclusive"
$ws.Range("N46").Value = "'2"

$ws.Range("A44").Value = "Giftmessage50letters"
$ws.Range("F44").Value = "Testing"
$ws.Range("G44").Value = "qa"
$ws.Range("AH44").Value = "By submitting this form and signing up for texts, you consent to receive marketing text messages (e.g. promos, cart reminders) from OXO at the number provided, including messages sent by autodialer."

$ws.Range("A23").Select()
$ws.Range("D36").Select()
